$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.490.85'
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").Value = '2.493.32'
$ws.Range("E3").Value = '  -0.47%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.20'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.99'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.53%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.511'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("E9").Value = '  +0.63%  '

$ws.Range("E10").Value = '  -0.59%  '

$ws.Range("E11").Value = '  -2.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.87'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.26%  '

$ws.Range("D13").Value = '2.951.95'
$ws.Range("E13").Value = '  -0.66%  '

$ws.Range("D14").Value = '69.393.07'
$ws.Range("E14").Value = '  +0.35%  '

$ws.Range("E15").Value = '  +0.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '24.16'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.31%  '

$ws.Range("D17").Value = '2.514.80'
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("E18").Value = '  -0.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.45'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '352.98'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("E22").Value = '  -3.90%  '

$ws.Range("E23").Value = '  -0.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.22'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.30%  '

$ws.Range("E25").Value = '  -2.67%  '

$ws.Range("D26").Value = '2.621.88'
$ws.Range("E26").Value = '  -1.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.65'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.38%  '

$ws.Range("E28").Value = '  +0.11%  '

$ws.Range("E29").Value = '  -1.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.54'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.56'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +134.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.19'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '438.51'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.59%  '

$ws.Range("E35").Value = '  -0.63%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.61%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '152.96'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.07'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.01%  '

$ws.Range("E39").Value = '  -1.39%  '

$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("E41").Value = '  -0.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.59'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.57'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.17'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.99%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.07'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '139.09'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.70%  '

$ws.Range("E47").Value = '  -0.92%  '

$ws.Range("E48").Value = '  -2.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0723'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.573'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0926'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.27%  '
